# Applies "Actualización automática 2025-07-02 17:15:08"
# - Adds julio (July) sales of 625.86 for DECORHOME S.C.C. (row 8)
# - Adds julio (July) sales of 6678.14 for PADILLA MIER BERTHA MARIETA (row 16)
# - Updates the derived totals / summary cells on both affected sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTA MENSUAL": raw monthly sales figures (column F = julio) ---
$wsVenta = $wb.Worksheets.Item("VENTA MENSUAL")

$wsVenta.Range("F8").Value = 625.86
$wsVenta.Range("F16").Value = 6678.14
$wsVenta.Range("F22").Value = 10245.14

# --- Sheet "VENTAS POR GRUPO": sales broken out by product group ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# PANELES PVC sale for DECORHOME S.C.C.
$wsGrupo.Range("Q8").Value = 625.86

# PORCELANATO sale for PADILLA MIER BERTHA MARIETA
$wsGrupo.Range("M16").Value = 6678.14

# Row 22 "x de 20" counters: the two clients above now each have one more
# product group populated, so the PORCELANATO (M22) and PANELES PVC (Q22)
# counters go from "1 de 20" to "2 de 20".
$wsGrupo.Range("M22").Value = "2 de 20"
$wsGrupo.Range("Q22").Value = "2 de 20"
